$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Text = "O"
